# Fix "Positive sign"->"string w/o leading spaces"
# Cell B3 on the "Cell data types" sheet currently holds the shared string
# "Positive Value". It should instead read "String without leading spaces".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "String without leading spaces"
